$wb = $excel.ActiveWorkbook

# The "TestResultExcelFilePath" output-path column (column H) is removed from
# the two NI-scenario sheets that drove the extra file-path input.
$wsProcess = $wb.Worksheets.Item("ProcessPayrollForNIWeekly")
$wsProcess.Columns("H").Delete()

$wsReports = $wb.Worksheets.Item("TestReports")
$wsReports.Columns("H").Delete()

# After the wide wrapped-text column is gone, the data rows no longer need
# to be as tall.
$wsProcess.Range("A3:A10").RowHeight = 30

# Update the cached cell selections left behind on each sheet.
[void]$wsProcess.Range("K4").Select()
[void]$wsReports.Range("I11").Select()

# Make "NI4WeeklyCat_B" the active sheet/tab, as in the saved workbook.
$ws2 = $wb.Worksheets.Item("NI4WeeklyCat_B")
$ws2.Activate()

$wb.Save()
